# Commit: "Fruta / hortaliza, semanal"
# A new weekly price record was inserted as row 35 on the "Hortaliza, Terminal
# Hortofrutícola Agro Chillán - Choclo" sheet, pushing the existing rows 35-111
# down to 36-112.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 35; this shifts rows 35..111 down
# to 36..112 (carrying all of their existing content/formatting with them).
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new record's data.
$ws.Range("A35").Value = 7
$ws.Range("B35").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C35").Value = "Ñuble"
$ws.Range("D35").Value = 44540
$ws.Range("E35").Value = 16
$ws.Range("F35").Value = 100112024
$ws.Range("G35").Value = "Choclo"
$ws.Range("H35").Value = "Choclero"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 450
$ws.Range("L35").Value = 500
$ws.Range("M35").Value = 475
$ws.Range("N35").Value = "$/unidad"
$ws.Range("O35").Value = "Región Metropolitana"
$ws.Range("P35").Value = 475
$ws.Range("Q35").Value = 1
$ws.Range("R35").Value = "Hortaliza"
